# "Adiciona teste com valores preenchidos a planilha de teste"
#
# Adds a new test case (TC010 - "Produto é cadastrado e modal é fechado")
# to the "Produto" sheet, right below the existing TC008/TC009 rows, and
# also tidies up a stray highlight style that was left on C10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Produto")

# Make room for the new test-case row: everything from row 12 down
# (the blank spacer row + the "Observações adicionais" block) shifts
# down by one row, carrying its formatting/merges with it.
$ws.Rows.Item(12).Insert()

# Seed the new row's formatting from the row above (TC009), which already
# has the right bordered-cell style for every column...
$ws.Range("A11:I11").Copy()
$ws.Range("A12:I12").PasteSpecial(-4122)   # xlPasteFormats

# ...except column G, which needs the date number format (as used by G4).
$ws.Range("G4").Copy()
$ws.Range("G12").PasteSpecial(-4122)       # xlPasteFormats

# Fill in the TC010 test case data.
$ws.Range("A12").Value = "TC010"
$ws.Range("B12").Value = 8
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = "Relógio de pulso"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 100
$ws.Range("G12").Value = 45625
$ws.Range("H12").Value = "Quantidade e data vazios e código, nome e valor preenchidos"
$ws.Range("I12").Value = "Produto é cadastrado e modal é fechado"

# C10 had a leftover white-fill/right-aligned highlight style; clear it
# back to the plain bordered look used by the rest of that row.
$ws.Range("E10").Copy()
$ws.Range("C10").PasteSpecial(-4122)       # xlPasteFormats
$ws.Range("C10").Value = 1
